$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Example - " prefix from the seven "Example - ..." headings.
#    ("Example" and " - " each occur exactly seven times in the document,
#    always together, so two global MatchCase replacements do the job
#    without touching the lowercase "example" occurring inside normal
#    sentences elsewhere in the document.)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Example", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute(" - ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from the last paragraph (the image /
#    changepic paragraph) to the start of the "Merging data with pages"
#    heading paragraph.
# ---------------------------------------------------------------------------
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13,[char]7) -eq "Merging data with pages") {
        $headingPara = $d.Paragraphs.Item($i)
        break
    }
}

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$bmRange = $headingPara.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 3) Rework the final paragraph: it currently holds the picture followed by
#    three runs of text spelling out the "[b.number;ope=changepic;...]" TBS
#    tag. That tag is moved into the picture's Description (wp:docPr/@descr,
#    i.e. InlineShape.AlternativeText) and gets an extra "tagpos=inside"
#    switch; the now-empty trailing text is removed, and a short italic
#    explanatory paragraph is appended after the picture paragraph.
# ---------------------------------------------------------------------------
$imgPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*changepic*") {
        $imgPara = $d.Paragraphs.Item($i)
        break
    }
}

$picture = $d.InlineShapes.Item($d.InlineShapes.Count)
$picture.AlternativeText = "[b.number;ope=changepic;from=pic_[val].png;tagpos=inside;adjust]"

$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pPrXml = "<w:pPr><w:rPr><w:noProof/><w:lang w:val=""en-US""/></w:rPr></w:pPr>"

$picRange = $picture.Range.Duplicate
$picXml = $picRange.WordOpenXML

$newParaXml = '<w:p ' + $xmlNs + '><w:pPr><w:rPr><w:i/><w:color w:val="800000"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:i/><w:color w:val="800000"/><w:lang w:val="en-US"/></w:rPr>' + `
              '<w:t>The image is merged using a TBS field which is placed in the Description property of the image.</w:t></w:r></w:p>'

$replacement = '<w:p ' + $xmlNs + '>' + $pPrXml + `
               '<w:r><w:rPr><w:noProof/><w:lang w:eastAsia="fr-FR"/></w:rPr><w:drawing>' + `
               '<wp:inline xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" distT="0" distB="0" distL="0" distR="0">' + `
               '<wp:extent cx="1132675" cy="1514475"/><wp:effectExtent l="0" t="0" r="0" b="0"/>' + `
               '<wp:docPr id="1" name="Image 1" descr="[b.number;ope=changepic;from=pic_[val].png;tagpos=inside;adjust]"/>' + `
               '<wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr>' + `
               '<a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
               '<pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/>' + `
               '<pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr>' + `
               '<pic:blipFill><a:blip r:embed="rId11" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + `
               '<a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst>' + `
               '</a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' + `
               '<pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="1132675" cy="1514475"/></a:xfrm>' + `
               '<a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic>' + `
               '</a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>' + `
               $newParaXml

$imgPara.Range.InsertXML($replacement)
